# Add the missing "Moneda" (currency) values in column D for rows 34-39,
# matching the style used by the rest of column D (centered alignment).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = 34..39
foreach ($r in $rows) {
    $cell = $ws.Range("D$r")
    $cell.Value = "USD"
    $cell.HorizontalAlignment = -4108  # xlCenter
}

# Update the saved view state: scroll position and active selection.
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E38").Select()
